$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "19.947.55"
$ws.Range("E2").Value = "  -8.34%  "
$ws.Range("D3").Value = "1.419.38"
$ws.Range("E3").Value = "  -7.79%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").Value = "'273.22"
$ws.Range("E6").Value = "  -5.87%  "
$ws.Range("D7").Value = "'0.3744"
$ws.Range("E7").Value = "  -3.56%  "
$ws.Range("D8").Value = "'0.3084"
$ws.Range("E8").Value = "  -3.33%  "
$ws.Range("D9").Value = "'39.77"
$ws.Range("E9").Value = "  -7.85%  "
$ws.Range("D10").Value = "'1.013"
$ws.Range("D11").Value = "'0.06596"
$ws.Range("E11").Value = "  -8.31%  "
$ws.Range("D12").Value = "'1.005"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").Value = "'5.391"
$ws.Range("E13").Value = "  -4.29%  "
$ws.Range("D14").Value = "'17.07"
$ws.Range("E14").Value = "  -8.02%  "
$ws.Range("D15").Value = "'6.159"
$ws.Range("E15").Value = "  -7.02%  "
$ws.Range("D16").Value = "1.425.32"
$ws.Range("E16").Value = "  -7.71%  "
$ws.Range("D17").Value = "'0.00001007"
$ws.Range("E17").Value = "  -8.94%  "
$ws.Range("D18").Value = "'0.05842"
$ws.Range("E18").Value = "  -11.11%  "
$ws.Range("D19").Value = "'74.97"
$ws.Range("E19").Value = "  -10.14%  "
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").Value = "'5.657"
$ws.Range("E21").Value = "  -7.88%  "
$ws.Range("D22").Value = "'14.45"
$ws.Range("E22").Value = "  -6.03%  "
$ws.Range("D23").Value = "'10.99"
$ws.Range("E23").Value = "  +0.74%  "
$ws.Range("D24").Value = "'2.335"
$ws.Range("E24").Value = "  -1.79%  "
$ws.Range("D25").Value = "19.973.94"
$ws.Range("E25").Value = "  -8.24%  "
$ws.Range("D26").Value = "'2.290"
$ws.Range("E26").Value = "  -4.11%  "
$ws.Range("D27").Value = "'139.18"
$ws.Range("E27").Value = "  -4.54%  "
$ws.Range("D28").Value = "'16.90"
$ws.Range("E28").Value = "  -8.00%  "
$ws.Range("D29").Value = "1.585.92"
$ws.Range("E29").Value = "  -7.72%  "
$ws.Range("D30").Value = "'109.15"
$ws.Range("E30").Value = "  -7.05%  "
$ws.Range("E31").Value = "  -20.12%  "
$ws.Range("D32").Value = "'0.9002"
$ws.Range("E32").Value = "  -6.99%  "
$ws.Range("D33").Value = "'5.432"
$ws.Range("E33").Value = "  -7.69%  "
$ws.Range("D34").Value = "'0.07774"
$ws.Range("E34").Value = "  -5.38%  "
$ws.Range("D35").Value = "'8.415"
$ws.Range("E35").Value = "  -5.84%  "
$ws.Range("D36").Value = "'11.31"
$ws.Range("E36").Value = "  +6.09%  "
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("D38").Value = "'4.763"
$ws.Range("E38").Value = "  -7.15%  "
$ws.Range("D39").Value = "'0.05675"
$ws.Range("E39").Value = "  -6.81%  "
$ws.Range("D40").Value = "'0.1916"
$ws.Range("E40").Value = "  -6.10%  "
$ws.Range("D41").Value = "'0.02023"
$ws.Range("E41").Value = "  -8.20%  "
$ws.Range("D42").Value = "'1.092"
$ws.Range("E42").Value = "  -8.02%  "
$ws.Range("D43").Value = "'1.258"
$ws.Range("E43").Value = "  -15.16%  "
$ws.Range("D44").Value = "'0.5325"
$ws.Range("E44").Value = "  -7.36%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'12.29"
$ws.Range("E45").Value = "  -5.58%  "
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "'3.533"
$ws.Range("E46").Value = "  -5.72%  "
$ws.Range("D47").Value = "'0.5133"
$ws.Range("E47").Value = "  -6.86%  "
$ws.Range("D48").Value = "'1.795"
$ws.Range("E48").Value = "  -3.76%  "
$ws.Range("D49").Value = "'109.49"
$ws.Range("E49").Value = "  -7.32%  "
$ws.Range("D50").Value = "'1.051"
$ws.Range("E50").Value = "  -7.93%  "
$ws.Range("D51").Value = "'1.002"
$ws.Range("E51").Value = "  +0.10%  "
